# Update the LR-pairs sheet (Igfbp4-Fzd8) with the recomputed NATMI values
# following Dr Hou's advice. Sending/Target clusters cycle through
# ECs / FAPs / sCs, and three new rows (sCs -> *) are appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data keyed by row number; columns A..T in order.
$data = @{}
$data[2]  = @("ECs",  "Igfbp4", "Fzd8", "ECs",  3, 1, 3.910463666666666, 11.731391,          0.02584512419166262, 0.02584512419166263, 2, 0.6666666666666666, 1.027778333333333, 3.083335,           0.08020467841353289, 0.08020467841353289, 4.019089829887221,  36.17180846898499,  0.00207289987435012,  0.00207289987435012)
$data[3]  = @("ECs",  "Igfbp4", "Fzd8", "FAPs", 3, 1, 3.910463666666666, 11.731391,          0.02584512419166262, 0.02584512419166263, 3, 1,                  7.273511666666667, 21.820535,          0.567602609669802,   0.567602609669802,   28.44280310157611,  255.985227914185,   0.01466975993842784,  0.01466975993842784)
$data[4]  = @("ECs",  "Igfbp4", "Fzd8", "sCs",  3, 1, 3.910463666666666, 11.731391,          0.02584512419166262, 0.02584512419166263, 3, 1,                  4.513153666666667, 13.539461,          0.3521927119166651,  0.3521927119166651,  17.64852343558344,  158.836710920251,   0.009102464378884668, 0.009102464378884666)
$data[5]  = @("FAPs", "Igfbp4", "Fzd8", "ECs",  3, 1, 125.0119883333333, 375.035965,         0.826232037766454,   0.8262320377664542,  2, 0.6666666666666666, 1.027778333333333, 3.083335,           0.08020467841353289, 0.08020467841353289, 128.4846130159195,  1156.361517143275,  0.06626767488401641,  0.06626767488401641)
$data[6]  = @("FAPs", "Igfbp4", "Fzd8", "FAPs", 3, 1, 125.0119883333333, 375.035965,         0.826232037766454,   0.8262320377664542,  3, 1,                  7.273511666666667, 21.820535,          0.567602609669802,   0.567602609669802,   909.2761556156973,  8183.485400541275,  0.4689714608290377,   0.4689714608290378)
$data[7]  = @("FAPs", "Igfbp4", "Fzd8", "sCs",  3, 1, 125.0119883333333, 375.035965,         0.826232037766454,   0.8262320377664542,  3, 1,                  4.513153666666667, 13.539461,          0.3521927119166651,  0.3521927119166651,  564.1983135238739,  5077.784821714865,  0.2909929020533999,   0.2909929020533999)
$data[8]  = @("sCs",  "Igfbp4", "Fzd8", "ECs",  3, 1, 22.381277,         67.14383099999999, 0.1479228380418832,  0.1479228380418833,  2, 0.6666666666666666, 1.027778333333333, 3.083335,           0.08020467841353289, 0.08020467841353289, 23.00299157293166,  207.026924156385,   0.01186410365516635,  0.01186410365516636)
$data[9]  = @("sCs",  "Igfbp4", "Fzd8", "FAPs", 3, 1, 22.381277,         67.14383099999999, 0.1479228380418832,  0.1479228380418833,  3, 1,                  7.273511666666667, 21.820535,          0.567602609669802,   0.567602609669802,   162.7904793743983,  1465.114314369585,  0.0839613889023364,   0.08396138890233641)
$data[10] = @("sCs",  "Igfbp4", "Fzd8", "sCs",  3, 1, 22.381277,         67.14383099999999, 0.1479228380418832,  0.1479228380418833,  3, 1,                  4.513153666666667, 13.539461,          0.3521927119166651,  0.3521927119166651,  101.0101423572323,  909.0912812150908,  0.0520973454843805,   0.0520973454843805)

foreach ($r in $data.Keys) {
    $rowVals = $data[$r]
    for ($i = 0; $i -lt $rowVals.Length; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $rowVals[$i]
    }
}
